$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 573.6667
$ws.Range("J17").Value = 575.96875
$ws.Range("L17").Value = 1727.90625
$ws.Range("N17").Value = -2063.90625

$ws.Range("H40").Value = 1775
$ws.Range("I40").Value = 2450
$ws.Range("K40").Value = 2450
$ws.Range("M40").Value = -2275

$ws.Range("H76").Value = 25002660
$ws.Range("I76").Value = 29414408
$ws.Range("J76").Value = 2746.6667
$ws.Range("K76").Value = 29414408
$ws.Range("L76").Value = 2746.6667
$ws.Range("M76").Value = -29414093
$ws.Range("N76").Value = -3376.6667

$ws.Range("H79").Value = 25002660
$ws.Range("I79").Value = 29414408
$ws.Range("J79").Value = 2746.6667
$ws.Range("K79").Value = 29414408
$ws.Range("L79").Value = 2746.6667
$ws.Range("M79").Value = -29413316
$ws.Range("N79").Value = -4930.6667

$ws.Range("H138").Value = 2323.5
$ws.Range("I138").Value = 1494.92
$ws.Range("J138").Value = 3309.9048
$ws.Range("K138").Value = 4484.76
$ws.Range("L138").Value = 9929.714399999999
$ws.Range("M138").Value = 655.2399999999998
$ws.Range("N138").Value = -20209.7144

$ws.Range("H141").Value = 783.5139
$ws.Range("I141").Value = 700.4179
$ws.Range("J141").Value = 1897
$ws.Range("K141").Value = 2101.2537
$ws.Range("L141").Value = 5691
$ws.Range("M141").Value = 3078.7463
$ws.Range("N141").Value = -16051

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1269.42
$ws.Range("I32").Value = 1240.5862
$ws.Range("J32").Value = 1462.3846
$ws.Range("K32").Value = 1240.5862
$ws.Range("L32").Value = 1462.3846
$ws.Range("M32").Value = -953.5862
$ws.Range("N32").Value = -2036.3846

$ws.Range("H110").Value = 2010.1578
$ws.Range("I110").Value = 956.7143
$ws.Range("J110").Value = 4959.8
$ws.Range("K110").Value = 956.7143
$ws.Range("L110").Value = 4959.8
$ws.Range("M110").Value = 1088.2857
$ws.Range("N110").Value = -9049.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 4500
$ws.Range("J8").Value = 4500
$ws.Range("L8").Value = 4500
$ws.Range("N8").Value = -4780

$ws.Range("H140").Value = 60000
$ws.Range("J140").Value = 60000
$ws.Range("L140").Value = 60000
$ws.Range("N140").Value = -70360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 12975
$ws.Range("J43").Value = 12975
$ws.Range("L43").Value = 12975
$ws.Range("N43").Value = -13343

$ws.Range("H101").Value = 12975
$ws.Range("J101").Value = 12975
$ws.Range("L101").Value = 12975
$ws.Range("N101").Value = -19465

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2297823.8
$ws.Range("I5").Value = 2959065.5
$ws.Range("K5").Value = 8877196.5
$ws.Range("M5").Value = -8877084.5

$ws.Range("H122").Value = 1240.4348
$ws.Range("I122").Value = 305.9091
$ws.Range("J122").Value = 2097.0833
$ws.Range("K122").Value = 2753.1819
$ws.Range("L122").Value = 18873.7497
$ws.Range("M122").Value = -303.1819
$ws.Range("N122").Value = -23773.7497

$ws.Range("H131").Value = 7828933.5
$ws.Range("I131").Value = 27833692
$ws.Range("J131").Value = 983.8913
$ws.Range("K131").Value = 83501076
$ws.Range("L131").Value = 2951.6739
$ws.Range("M131").Value = -83496036
$ws.Range("N131").Value = -13031.6739

$ws.Range("H132").Value = 2439.875
$ws.Range("I132").Value = 2293.2727
$ws.Range("J132").Value = 2563.923
$ws.Range("K132").Value = 20639.4543
$ws.Range("L132").Value = 23075.307
$ws.Range("M132").Value = -18109.4543
$ws.Range("N132").Value = -28135.307

$ws.Range("H135").Value = 2297823.8
$ws.Range("I135").Value = 2959065.5
$ws.Range("K135").Value = 26631589.5
$ws.Range("M135").Value = -26629054.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 1993.3572
$ws.Range("I43").Value = 623.61536
$ws.Range("J43").Value = 19800
$ws.Range("K43").Value = 623.61536
$ws.Range("L43").Value = 19800
$ws.Range("M43").Value = -472.61536
$ws.Range("N43").Value = -20102

$ws.Range("H46").Value = 7760
$ws.Range("I46").Value = 4013.6667
$ws.Range("J46").Value = 18999
$ws.Range("K46").Value = 4013.6667
$ws.Range("L46").Value = 18999
$ws.Range("M46").Value = -3857.6667
$ws.Range("N46").Value = -19311

$ws.Range("H57").Value = 26392.5
$ws.Range("J57").Value = 26392.5
$ws.Range("L57").Value = 26392.5
$ws.Range("N57").Value = -28032.5

$ws.Range("H80").Value = 13879.579
$ws.Range("I80").Value = 5172.727
$ws.Range("J80").Value = 25851.5
$ws.Range("K80").Value = 5172.727
$ws.Range("L80").Value = 25851.5
$ws.Range("M80").Value = -4174.727
$ws.Range("N80").Value = -27847.5

$ws.Range("H83").Value = 13879.579
$ws.Range("I83").Value = 5172.727
$ws.Range("J83").Value = 25851.5
$ws.Range("K83").Value = 25863.635
$ws.Range("L83").Value = 129257.5
$ws.Range("M83").Value = -20871.635
$ws.Range("N83").Value = -139241.5

$ws.Range("H126").Value = 16216.714
$ws.Range("I126").Value = 18586.166
$ws.Range("K126").Value = 55758.49800000001
$ws.Range("M126").Value = -53288.49800000001

$ws.Range("H132").Value = 19914506
$ws.Range("I132").Value = 19048710
$ws.Range("J132").Value = 22728338
$ws.Range("K132").Value = 57146130
$ws.Range("L132").Value = 68185014
$ws.Range("M132").Value = -57143600
$ws.Range("N132").Value = -68190074

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 943.9
$ws.Range("I46").Value = 933.3333
$ws.Range("J46").Value = 948.4286
$ws.Range("K46").Value = 933.3333
$ws.Range("L46").Value = 948.4286
$ws.Range("M46").Value = -745.3333
$ws.Range("N46").Value = -1324.4286

$ws.Range("H136").Value = 3004287.8
$ws.Range("I136").Value = 3087601.5
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 9262804.5
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -9260254.5
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 14031398
$ws.Range("I126").Value = 21826064
$ws.Range("J126").Value = 998.8
$ws.Range("K126").Value = 65478192
$ws.Range("L126").Value = 2996.4
$ws.Range("M126").Value = -65475722
$ws.Range("N126").Value = -7936.4

$ws.Range("H128").Value = 53975
$ws.Range("J128").Value = 53975
$ws.Range("L128").Value = 53975
$ws.Range("N128").Value = -63935

$ws.Range("H136").Value = 3326.5
$ws.Range("I136").Value = 916.8929000000001
$ws.Range("J136").Value = 6137.7085
$ws.Range("K136").Value = 2750.6787
$ws.Range("L136").Value = 18413.1255
$ws.Range("M136").Value = -200.6787000000004
$ws.Range("N136").Value = -23513.1255
